# Add season stats for the new M3 (2021) seasons, and correct a value in
# the previously-last row (season 11 / row 13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the num_matches value for the existing last row (season 11 / row 13)
$ws.Range("E13").Value = 1182353

# New season rows (14-20), same shape as the existing data rows.
$newRows = @(
    @{ Row = 14; A = 12; B = "M3_01 Wolf 2021";    C = 9637; D = 10653; E = 808651; F = 9916; G = 10044; H = 10295 },
    @{ Row = 15; A = 13; B = "M3_02 Love 2021";    C = 9684; D = 10714; E = 917027; F = 9975; G = 10097; H = 10325 },
    @{ Row = 16; A = 14; B = "M3_03 Bear 2021";    C = 9637; D = 10576; E = 766502; F = 9914; G = 10026; H = 10230 },
    @{ Row = 17; A = 15; B = "M3_04 Elf 2021";     C = 9686; D = 10678; E = 944323; F = 9992; G = 10102; H = 10323 },
    @{ Row = 18; A = 16; B = "M3_05 Viper 2021";   C = 9701; D = 10753; E = 956484; F = 9998; G = 10106; H = 10300 },
    @{ Row = 19; A = 17; B = "M3_06 Magic 2021";   C = 9681; D = 10632; E = 869262; F = 9974; G = 10082; H = 10278 },
    @{ Row = 20; A = 18; B = "M3_07 Griffin 2021"; C = 9669; D = 10633; E = 856695; F = 9958; G = 10067; H = 10287 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
}

# Replicate column-A styling (bold/bordered/centered) from the last existing
# row onto the newly added rows, without introducing new style entries.
$ws.Range("A13").Copy()
$ws.Range("A14:A20").PasteSpecial(-4122)
